# Add a new "503A" indicator column (L) to the facility_info sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("facility_info")

# Header + 64 data rows for the new column L ("503A" Y/NA flag).
$values = @(
    "503A",
    "Y","NA","Y","NA","Y","NA","NA","NA","NA",
    "NA","NA","NA","NA","NA","Y","NA","NA","NA","NA",
    "NA","NA","NA","NA","NA","NA","NA","NA","NA","NA",
    "NA","NA","NA","NA","NA","NA","NA","NA","NA","NA",
    "NA","NA","NA","Y","NA","NA","NA","NA","NA","NA",
    "NA","NA","Y","NA","NA","NA","NA","NA","NA","NA",
    "NA","NA","Y","NA"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 12).Value = $values[$i]
}

# Re-establish the AutoFilter over the now-wider A1:L64 range.
$ws.AutoFilterMode = $false
$ws.Range("A1:L64").AutoFilter() | Out-Null

# Keep the workbook's hidden _FilterDatabase defined name in sync.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=facility_info!`$A`$1:`$L`$64"
    }
}

# Mirror the author's new selection on the sheet.
$ws.Range("L65").Select() | Out-Null
